# [GPR] Ata de reunião - Retrospective
# Alteração da data da reunião de review: 01/06/2015 -> 03/06/2015.

$d = $word.ActiveDocument

# Locate the meeting date cell's text ("01/06/2015") in the document.
$dateRng = $d.Content
$dateRng.Find.Execute("01/06/2015", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)

if ($dateRng.Find.Found) {
    $dayStart = $dateRng.Start
    $dayEnd = $dayStart + 2

    # Replace just the "day" portion of the date ("01" -> "03"); this is
    # the only part that changed, and it is where the user's cursor/last
    # edit lands.
    $dayRng = $d.Range($dayStart, $dayEnd)
    $dayRng.Text = "03"

    # Word tracks the location of the most recent edit with the special,
    # single-instance "_GoBack" bookmark. Adding it again simply moves it
    # (rather than duplicating it), which both places it right after the
    # edited "03" and removes it from wherever it previously sat in the
    # document.
    $markRng = $d.Range($dayStart + 2, $dayStart + 2)
    $d.Bookmarks.Add("_GoBack", $markRng)
}
